$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previously-existing data rows (A2:F20) first so stale values
# from the old (shorter) table can't linger once we rewrite the sheet with
# the refreshed/extended quarterly series below.
$ws.Range("A2:F20").Clear()

# Quarters from 2018Q4 through 2025Q2 (27 data rows total), replacing /
# extending the previous 2020Q2..2024Q4 (19 rows) recurrence data set.
$quarters = @(
    "2018Q4","2019Q1","2019Q2","2019Q3","2019Q4","2020Q1",
    "2020Q2","2020Q3","2020Q4","2021Q1","2021Q2","2021Q3","2021Q4",
    "2022Q1","2022Q2","2022Q3","2022Q4",
    "2023Q1","2023Q2","2023Q3","2023Q4",
    "2024Q1","2024Q2","2024Q3","2024Q4",
    "2025Q1","2025Q2"
)

$data = @(
    @(524, 36, 488, 6.94980694980695),
    @(407, 29, 378, 5.534351145038168),
    @(321, 17, 304, 4.176904176904177),
    @(568, 20, 548, 6.230529595015576),
    @(1032, 23, 1009, 4.049295774647888),
    @(727, 51, 676, 4.941860465116279),
    @(658, 21, 637, 2.888583218707015),
    @(897, 16, 881, 2.43161094224924),
    @(673, 50, 623, 5.574136008918618),
    @(525, 21, 504, 3.12035661218425),
    @(713, 36, 677, 6.857142857142858),
    @(700, 40, 660, 5.610098176718092),
    @(747, 47, 700, 6.714285714285714),
    @(735, 49, 686, 6.559571619812584),
    @(752, 67, 685, 9.115646258503402),
    @(821, 77, 744, 10.23936170212766),
    @(726, 58, 668, 7.064555420219245),
    @(797, 51, 746, 7.024793388429752),
    @(763, 58, 705, 7.277289836888332),
    @(748, 51, 697, 6.684141546526867),
    @(672, 40, 632, 5.347593582887701),
    @(653, 50, 603, 7.440476190476191),
    @(596, 62, 534, 9.494640122511486),
    @(715, 64, 651, 10.73825503355705),
    @(738, 49, 689, 6.853146853146853),
    @(644, 63, 581, 8.536585365853659),
    @(63, 19, 44, 2.950310559006211)
)

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $row = $i + 2
    $q = $quarters[$i]
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $q
    $ws.Cells.Item($row, 2).Value = $q
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}

$wb.Save()
